$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old training-metric columns (E:H); this shifts test->E:H, vs->I:L, Kernel->M
$ws.Range("E1:H1").EntireColumn.Delete()

# Update header labels for the shifted test/vs columns
$ws.Range("E1").Value = "r2_test"
$ws.Range("F1").Value = "mse_test"
$ws.Range("G1").Value = "mape_test"
$ws.Range("H1").Value = "rmse_test"
$ws.Range("I1").Value = "r2_vs"
$ws.Range("J1").Value = "mse_vs"
$ws.Range("K1").Value = "mape_vs"
$ws.Range("L1").Value = "rmse_vs"
$ws.Range("M1").Value = "Kernel"

# Row 2
$ws.Range("A2").Value = 0.8279410895582564
$ws.Range("B2").Value = 1.150558726382374
$ws.Range("C2").Value = 1.928099377909652
$ws.Range("D2").Value = 1.072641005361241
$ws.Range("E2").Value = 0.812596694487752
$ws.Range("F2").Value = 0.6703307575493209
$ws.Range("G2").Value = 0.3912529869704837
$ws.Range("H2").Value = 0.8187372945880266
$ws.Range("I2").Value = 0.8381668518268482
$ws.Range("J2").Value = 1.004372750444871
$ws.Range("K2").Value = 1.023605295675912
$ws.Range("L2").Value = 1.002183990315586
$ws.Range("M2").Value = "RBF"

# Row 3
$ws.Range("A3").Value = 0.7543574647673595
$ws.Range("B3").Value = 1.642612764180538
$ws.Range("C3").Value = 1.356114230960144
$ws.Range("D3").Value = 1.281644554539416
$ws.Range("E3").Value = 0.7675885757354772
$ws.Range("F3").Value = 0.8313221886055374
$ws.Range("G3").Value = 3.478574488141885
$ws.Range("H3").Value = 0.9117687144257239
$ws.Range("I3").Value = 0.7667261123846064
$ws.Range("J3").Value = 1.447749974316511
$ws.Range("K3").Value = 0.8687693026741286
$ws.Range("L3").Value = 1.203224822847547
$ws.Range("M3").Value = "Matern_0.5"

# Row 4
$ws.Range("A4").Value = 0.8279408056191435
$ws.Range("B4").Value = 1.150560625084529
$ws.Range("C4").Value = 1.928096498863575
$ws.Range("D4").Value = 1.072641890420344
$ws.Range("E4").Value = 0.8125967518320625
$ws.Range("F4").Value = 0.6703305524320474
$ws.Range("G4").Value = 0.3912274299785801
$ws.Range("H4").Value = 0.818737169323616
$ws.Range("I4").Value = 0.8381666317317387
$ws.Range("J4").Value = 1.004374116404398
$ws.Range("K4").Value = 1.023603835096794
$ws.Range("L4").Value = 1.002184671806748
$ws.Range("M4").Value = "RationalQuadratic"

# Row 5
$ws.Range("A5").Value = 0.8279410670123399
$ws.Range("B5").Value = 1.15055887714702
$ws.Range("C5").Value = 1.928099017629703
$ws.Range("D5").Value = 1.072641075638548
$ws.Range("E5").Value = 0.8125967065228952
$ws.Range("F5").Value = 0.6703307145003117
$ws.Range("G5").Value = 0.3912504615958605
$ws.Range("H5").Value = 0.8187372682981469
$ws.Range("I5").Value = 0.838166835407087
$ws.Range("J5").Value = 1.004372852349585
$ws.Range("K5").Value = 1.023605121471817
$ws.Range("L5").Value = 1.002184041156905
$ws.Range("M5").Value = "ExpSineSquared"

# Row 6
$ws.Range("A6").Value = 0.281932692313457
$ws.Range("B6").Value = 4.801719392895841
$ws.Range("C6").Value = 2.625724923012809
$ws.Range("D6").Value = 2.191282590834838
$ws.Range("E6").Value = 0.2983028613606485
$ws.Range("F6").Value = 2.509929978175148
$ws.Range("G6").Value = 8.835420966828657
$ws.Range("H6").Value = 1.584275852929391
$ws.Range("I6").Value = 0.296803636346945
$ws.Range("J6").Value = 4.364194071719984
$ws.Range("K6").Value = 4.498193682469465
$ws.Range("L6").Value = 2.0890653584127
$ws.Range("M6").Value = "DotProduct"

